# Slide 1, subtitle placeholder (Shape 31 / p:ph type="subTitle"):
#   - "Basics, Services and Random Topics" is split into two runs
#     ("Basics, Services and Random " + "Topics") within the first paragraph.
#   - A new second paragraph "Shreyansh Jain" is added, itself split into
#     two runs ("Shreyansh" + " Jain") with lang="en-US".
#
# Build order is deliberately chosen so each chunk of text picks up the
# correct language id: the "Shreyansh"/" Jain" paragraph is typed first
# (as the sole run) so LanguageID="en-US" lands on it, then the title
# paragraph text is inserted *before* it (inheriting that run's
# properties at first, which is why LanguageID="en" is re-applied to fix
# it back to "en" afterwards), and finally the title text is split into
# its two runs.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# 1) Seed the text range with the new second paragraph's first word and
#    mark it (and anything that inherits from it) as English (US).
$tr.Text = "Shreyansh"
$tr.LanguageID = "en-US"
[void]$tr.InsertAfter(" Jain")

# 2) Prepend the (unchanged) title text plus a paragraph break, ahead of
#    the "Shreyansh Jain" paragraph just built.
[void]$tr.InsertBefore("Basics, Services and Random Topics" + [char]13)

# 3) The just-inserted title paragraph inherited "en-US" from the run it
#    was placed before; restore it to "en".
$tr.LanguageID = "en"

# 4) Split the title paragraph's single run into two runs, matching the
#    "Basics, Services and Random " / "Topics" boundary.
$full = $sh.TextFrame.TextRange
$titleFirstPart = $full.Characters(1, 28)
$titleFirstPart.Text = "Basics, Services and Random "
